$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits -------------------------------------------------
# DepartmentName for row 2 (E2): "school" -> "Store"
$ws.Range("E2").Value = "Store"

# LocationOfficeName for row 2 (G2): "Ghana" -> "Dhaka Office"
$ws.Range("G2").Value = "Dhaka Office"

# Email for row 5 (I5): "employee4@mail.cm" -> "shumonsb@gmail.com"
$ws.Range("I5").Value = "shumonsb@gmail.com"

# --- Hyperlinks ---------------------------------------------------------
# Only I5's hyperlink display text should change (its target mailto:
# address stays the same: employee4@mail.cm). The other hyperlinks on the
# sheet must keep pointing at the same addresses/display text.
# This host's Hyperlinks collection only reliably supports bulk delete +
# re-add (per-item mutation/delete does not persist), so rebuild the full
# set of hyperlinks, changing only I5's display text.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:shumonsb@gmail.com", "", "", "shumonsb@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:employee@gmail.com", "", "", "employee@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:example@email.com", "", "", "example@email.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:mono@gmail.com", "", "", "mono@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:employee3@mail.com", "", "", "employee3@mail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:sb2@gmail.com", "", "", "sb2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I5"), "mailto:employee4@mail.cm", "", "", "shumonsb@gmail.com")

# --- Selection / active cell --------------------------------------------
$ws.Range("I11").Select()
